{"js": "// The document starts as a single paragraph holding the run \"test\"\n// (with an en-US language tag). The edit:\n//   1. Replaces that paragraph's text with a run-split Russian greeting,\n//      interleaved with the w:proofErr spell/grammar-check markers Word\n//      leaves behind after its proofing pass (\"\u0421\u0430\u043b\u0430\u043c\u0430\u043b\u0435\u0439\u043a\u0443\u043c , \u044d\u0442\u043e\n//      \u041f\u0430\u0445\u0440\u0443\u0434\u0438\u043d\"), dropping the old language formatting.\n//   2. Splits the trailing bookmark (\"_GoBack\") out into its own new,\n//      empty paragraph right after it.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\nconst range = firstParagraph.getRange(\"Whole\");\n\n// Office.js insertOoxml requires a full flat-OPC (\"pkg:package\") payload;\n// build one whose /word/document.xml part supplies the replacement\n// paragraphs (the run text + w:proofErr markers, then a second paragraph\n// carrying the relocated bookmark). Replacing the paragraph's whole range\n// (which includes its end-of-paragraph mark) with two paragraphs' worth\n// of content turns the one paragraph into two.\nconst newParagraphsXml =\n  \"<w:p>\" +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    \"<w:r><w:t>\\u0421\\u0430\\u043B\\u0430\\u043C\\u0430\\u043B\\u0435\\u0439\\u043A\\u0443\\u043C</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> ,</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> \\u044D\\u0442\\u043E </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    \"<w:r><w:t>\\u041F\\u0430\\u0445\\u0440\\u0443\\u0434\\u0438\\u043D</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n  \"</w:p>\" +\n  \"<w:p>\" +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n  \"</w:p>\";\n\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      \"<pkg:xmlData>\" +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          newParagraphsXml +\n        \"</w:document>\" +\n      \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n  \"</pkg:package>\";\n\nrange.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The first (and only) paragraph currently holds \"test\". Replace its whole\n# range (text + paragraph mark) with the new paragraph -- run-split text\n# with spelling/grammar proofing marks -- followed by a second, new\n# paragraph that carries the bookmark that used to sit at the end of the\n# original paragraph.\n$xml = @'\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:p><w:proofErr w:type=\"spellStart\"/><w:proofErr w:type=\"gramStart\"/><w:r><w:t>\u0421\u0430\u043b\u0430\u043c\u0430\u043b\u0435\u0439\u043a\u0443\u043c</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> ,</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\"> \u044d\u0442\u043e </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>\u041f\u0430\u0445\u0440\u0443\u0434\u0438\u043d</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p><w:p><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p></w:document></pkg:xmlData></pkg:part></pkg:package>\n'@\n\n$r = $d.Paragraphs(1).Range\n$r.InsertXML($xml)\n"}
